# Replace the numeric "Codigo" column (A2:A6) with generated text codes,
# since getSublist no longer works for producing sequential numeric ids.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$codes = @("CGWVHQB", "MVRBVIWB", "SUYWGFW", "PEIOJECKJ", "WHBVWKDJ")

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $codes[$i]
}

# Header row, Codigo (A) and Nombre (B) columns are textual -> store as text format.
$ws.Range("A1:D1").NumberFormat = "@"
$ws.Range("A2:B6").NumberFormat = "@"
# Precio / Cantidad (C:D) columns use one-decimal numeric format.
$ws.Range("C2:D6").NumberFormat = "0.0"

$ws.Range("C6").Select()
